## loop fix for child entries
## The "child loop" that stamps out the Schema/Field/Comment formulas for
## each Source Topic stopped one row short (row 14 instead of row 15), and
## the next row (16) was left with hard-coded literal values instead of the
## same formulas. This brings rows 15/16 in line with the rest of the
## Schema block (rows 3-14):
##   - Row 15 ("Zenga") gets the same live formulas as the other child rows.
##   - Row 16 ("Zenga" / "RollaZenga" / "RollaZenga File") is filled in with
##     the literal values that loop iteration produced.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 15 - continues the pattern from rows 3-14: Schema name in I, and the
# J/N columns driven off the same formulas used by the rows above.
$ws.Range("I15").Value = "Zenga"
$ws.Range("J15").Formula = "=I15&""""&""name"""
$ws.Range("N15").Formula = "=""This Is ""&I15&"" ""&J15&"" field"""

# Row 16 - the loop's final (previously-skipped) child entry, written out
# with literal values.
$ws.Range("I16").Value = "Zenga"
$ws.Range("J16").Value = "RollaZenga"
$ws.Range("N16").Value = "RollaZenga File"

# Leave the selection where the edit finished, matching the saved view state.
$ws.Range("I16").Select()
